$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 22-24 down to 120-122 by inserting 98 blank rows before row 22
$ws.Range("A22:A119").EntireRow.Insert()

# Fill new data rows 21-45, column A first (becomes shared-string indices 44-68)
$ws.Range("A21").Value = "He tried again and once more sliced out of bounds"
$ws.Range("A22").Value = "Somebody brought a light"
$ws.Range("A23").Value = "The prediction was correct"
$ws.Range("A24").Value = "He sized up the situation and shook his head"
$ws.Range("A25").Value = "The rancher was trembling"
$ws.Range("A26").Value = "The vision became even stronger now"
$ws.Range("A27").Value = "He felt a little sick at his stomach"
$ws.Range("A28").Value = "His heart beat faster"
$ws.Range("A29").Value = "I never heard that"
$ws.Range("A30").Value = "Come inside now"
$ws.Range("A31").Value = "It was my fault"
$ws.Range("A32").Value = "He scarcely saw them"
$ws.Range("A33").Value = "Then he smiled shyly"
$ws.Range("A34").Value = "His room will be ready shortly"
$ws.Range("A35").Value = "The ball broke up in confusion"
$ws.Range("A36").Value = "Such a pitiful end"
$ws.Range("A37").Value = "That was the day it ended"
$ws.Range("A38").Value = "The girl smiled and started forward"
$ws.Range("A39").Value = "That was gonna be fun collecting"
$ws.Range("A40").Value = "He cleared his throat"
$ws.Range("A41").Value = "I been riding train for a ways now"
$ws.Range("A42").Value = "This was a slightly different matter"
$ws.Range("A43").Value = "Now forget all this other"
$ws.Range("A44").Value = "Youth obeyed when commanded"
$ws.Range("A45").Value = "But the past was dead here as the present was dead"

# then column B (becomes shared-string indices 69-93)
$ws.Range("B21").Value = "If we wait until children are in junior high or high school , we will never manage it"
$ws.Range("B22").Value = "Improvement can be measured by the lessening distance between toes and head"
$ws.Range("B23").Value = "He places the hands on either side of the head , keeping the chin down on the chest"
$ws.Range("B24").Value = "He then pushes his seat into the air and the teacher guides it over"
$ws.Range("B25").Value = "It is very important for parents to understand that early training is imperative"
$ws.Range("B26").Value = "This stain often disrupts the normal cell activity or else colors only the outside"
$ws.Range("B27").Value = "A balanced resistance bridge and a pen recorder are all the electronic instrumentation needed"
$ws.Range("B28").Value = "The transducer is coupled to the body through a water bath , not shown"
$ws.Range("B29").Value = "His earlier love for literature and history remained with him for his entire life"
$ws.Range("B30").Value = "He proposed a fresh theory of alkalis which later was accepted in chemical practices"
$ws.Range("B31").Value = "The form of galvanic activity is halfway between the magnetic form and the electrical form"
$ws.Range("B32").Value = "He devised a detonating fuse in which a short wire was caused to glow by an electric current"
$ws.Range("B33").Value = "Cows receiving drug may not be officially tested under breed registry testing programs"
$ws.Range("B34").Value = "Several materials or combinations of materials can be used to construct a satisfactory feed bunk"
$ws.Range("B35").Value = "Here are some key areas to examine to make sure your pricing strategy will be on target"
$ws.Range("B36").Value = "This problem can force a change in marketing approach in many kinds of businesses"
$ws.Range("B37").Value = "No manufacturer has taken the initiative in pointing out the costs involved"
$ws.Range("B38").Value = "Computers are being used to keep branch inventories at more workable levels"
$ws.Range("B39").Value = "There may be possible economies at any one of a number of links in your marketing and distribution chain"
$ws.Range("B40").Value = "Are there individuals in your organization who can shepherd a new product through to commercialization"
$ws.Range("B41").Value = "Most marketing people agree it is going to take redoubled efforts to satisfy future requirements"
$ws.Range("B42").Value = "Every single problem touched on thus far is related to good marketing planning"
$ws.Range("B43").Value = "Beyond any question of curriculum and approach to subject must be the quality of the teachers themselves"
$ws.Range("B44").Value = "Only a few years ago a middle western college circulated a request for a teacher of interior design"
$ws.Range("B45").Value = "One solution is the aquisition of degrees in education but it is a poor substitute"

$ws.Range("A52").Select()
